$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values are plain decimals (e.g. "0.556") which
# Excel would otherwise auto-convert to numeric cells. The source data keeps
# every Price/Volume cell as text, so force those specific cells to Text
# format before writing them, then restore their style afterwards.
$textForceRefs = @("D5", "D6", "D8", "D9", "D10", "D11", "D14", "D15", "D17", "D18", "D20", "D21", "D23", "D24", "D26", "D27", "D30", "D31", "D32", "D33", "D37", "D40", "D43", "D44", "D49")
foreach ($r in $textForceRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.298.18"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.802.44"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "227.16"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "0.556"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "33.29"
$ws.Range("E8").Value = "  +3.09%  "
$ws.Range("D9").Value = "0.296"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "0.0689"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "0.0947"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.058.60"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.804.38"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "11.08"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "0.634"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "34.298.24"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "4.28"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "68.48"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = ("0.0" + [string][char]0x2083 + "0795")
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").Value = "244.24"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").Value = "11.29"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "4.17"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "167.23"
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").Value = "7.31"
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("D27").Value = "16.50"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "3.97"
$ws.Range("E30").Value = "  +5.94%  "
$ws.Range("D31").Value = "0.0526"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "1.24"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "3.79"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "1.403.56"
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("D37").Value = "0.669"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").Value = "84.83"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("E41").Value = "  +4.25%  "
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").Value = "0.938"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "13.94"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  +3.16%  "
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").Value = "1.958.48"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").Value = "104.98"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E51").Value = "  -0.62%  "

foreach ($r in $textForceRefs) {
    $ws.Range($r).Style = "Normal"
}
